$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7811
$ws.Range("F3").Value = 7811
$ws.Range("F5").Value = 7958
$ws.Range("F9").Value = 6860
$ws.Range("F10").Value = 3426
$ws.Range("F12").Value = 3752
$ws.Range("F13").Value = 53
$ws.Range("F15").Value = 53
$ws.Range("F17").Value = 91
$ws.Range("F18").Value = 485
$ws.Range("F20").Value = 76
$ws.Range("F25").Value = 3933
$ws.Range("F27").Value = 385
$ws.Range("F29").Value = 506
$ws.Range("F30").Value = 1547
$ws.Range("F31").Value = 85
$ws.Range("F32").Value = 72
$ws.Range("F33").Value = 2842
$ws.Range("F34").Value = 2011
$ws.Range("F35").Value = 44
$ws.Range("F37").Value = 70
$ws.Range("F38").Value = 104
$ws.Range("F39").Value = 3877
$ws.Range("F40").Value = 356
$ws.Range("F41").Value = 294
$ws.Range("F42").Value = 45
$ws.Range("F44").Value = 629
$ws.Range("F45").Value = 54
$ws.Range("F46").Value = 1494
$ws.Range("F48").Value = 574
$ws.Range("F49").Value = 667
$ws.Range("F50").Value = 14

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 46
$ws.Range("F8").Value = 48
$ws.Range("F18").Value = 326

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("E2").Value = "2024.05.25 00:00-06.30 23:59"
$ws.Range("F2").Value = 147

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("E2").Value = "2024.05.25 00:00-06.30 23:59"
$ws.Range("F2").Value = 147
$ws.Range("F4").Value = 46
$ws.Range("F5").Value = 7811
$ws.Range("F6").Value = 7811
$ws.Range("F7").Value = 7958
$ws.Range("F11").Value = 6860
$ws.Range("F12").Value = 3426
$ws.Range("F13").Value = 3752
$ws.Range("F15").Value = 53
$ws.Range("F16").Value = 91
$ws.Range("F17").Value = 485
$ws.Range("F19").Value = 48
$ws.Range("F20").Value = 76
$ws.Range("F26").Value = 3933
$ws.Range("F30").Value = 385
$ws.Range("F31").Value = 506
$ws.Range("F32").Value = 1547
$ws.Range("F33").Value = 85
$ws.Range("F34").Value = 72
$ws.Range("F35").Value = 2011
$ws.Range("F36").Value = 44
$ws.Range("F38").Value = 104
$ws.Range("F40").Value = 3877
$ws.Range("F41").Value = 356
$ws.Range("F42").Value = 294
$ws.Range("F44").Value = 45
$ws.Range("F46").Value = 54
$ws.Range("F47").Value = 1494
$ws.Range("F50").Value = 667
